$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.272.83'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '2.942.22'
$ws.Range("E3").Value = '  -2.64%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.73'
$ws.Range("E5").Value = '  -3.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.00'
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = '2.938.51'
$ws.Range("E9").Value = '  -2.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.72'
$ws.Range("E10").Value = '  -3.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.150'
$ws.Range("E11").Value = '  -3.47%  '
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.33'
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '65.248.97'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").Value = '3.431.85'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = '2.942.92'
$ws.Range("E19").Value = '  -2.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.82'
$ws.Range("E20").Value = '  +7.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '445.43'
$ws.Range("E21").Value = '  -3.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.686'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.24'
$ws.Range("E23").Value = '  -2.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.08'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.09'
$ws.Range("E26").Value = '  -4.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  -7.08%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.03'
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.40'
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.57'
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.10'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.111'
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  -2.41%  '
$ws.Range("E37").Value = '  -1.63%  '
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '44.43'
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("E40").Value = '  -9.71%  '
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("E42").Value = '  -7.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.298'
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.43'
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '386.10'
$ws.Range("E45").Value = '  -1.94%  '
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").Value = '2.706.65'
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.03'
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.19'
$ws.Range("E50").Value = '  +4.65%  '
$ws.Range("E51").Value = '  -0.50%  '
